# Store_template.xlsx edit: populate store21 record with real contact info,
# work schedule (with hyperlink email + time formatting), adjust layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 values -----------------------------------------------------
# A: storename (unchanged, still "store21")
$ws.Range("A2").Value = "store21"

# E: email -> also becomes a mailto hyperlink (adds Hyperlink style/font)
$ws.Range("E2").Value = "store21@gmail.com"
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:store21@gmail.com") | Out-Null

# F/G: telephone / mobile
$ws.Range("F2").Value = "066/44564545"
$ws.Range("G2").Value = "035/44565454"

# H: comment
$ws.Range("H2").Value = "No comment"

# B/C/D: street / zipcode / place -> blank (single space)
$ws.Range("B2").Value = " "
$ws.Range("C2").Value = " "
$ws.Range("D2").Value = " "

# I/J: start_work / end_work -> long JS-style date strings, wrapped + time numfmt
$ws.Range("I2").Value = "Tue Aug 13 2019 07:00:00 GMT+0200 (Central European Summer Time)"
$ws.Range("J2").Value = "Tue Aug 13 2019 14:00:00 GMT+0200 (Central European Summer Time)"
$ws.Range("I2").NumberFormat = "h:mm"
$ws.Range("I2").WrapText = $true
$ws.Range("J2").NumberFormat = "h:mm"
$ws.Range("J2").WrapText = $true

# K/L: time_duration / time_therapy
$ws.Range("K2").Value = 45
$ws.Range("L2").Value = 15

# M: superadmin stays 85
$ws.Range("M2").Value = 85

# --- Row height for the now-wrapped row --------------------------------
$ws.Rows.Item(2).RowHeight = 90

# --- Column width tweaks -------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 24.6    # column E -> ~25.42578125 "chars"
$ws.Columns.Item(9).ColumnWidth = 13.76   # column I -> ~14.7109375 "chars"

# --- Selection / view: drop the stale topLeftCell / old selection -------
$ws.Range("D2").Select() | Out-Null
